# Cleaned up notebook folders:
#  - Reposition the "kreativKOPF" title textbox slightly lower on the
#    cover slide.
#  - Remove the now-unused "Marketing Intelligence by kreativbox.io &
#    machinemind.io" subtitle textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point conversion (PowerPoint COM measures Top/Left in points,
# the OOXML stores EMUs; 914400 EMU per inch, 72 points per inch).
$emuPerPoint = 914400 / 72

$titleBox = $s.Shapes.Item("TextBox 3")
$titleBox.Top = 2105561 / $emuPerPoint

$subtitleBox = $s.Shapes.Item("TextBox 4")
$subtitleBox.Delete()
